# Update "Top 50 Cryptocurrencies" sheet with the latest market data
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Top 50 Cryptocurrencies")
$ws2 = $wb.Worksheets.Item("Top 5 by Market Cap")
$ws3 = $wb.Worksheets.Item("Summary")

# Each entry: Row, Current Price (USD), Market Capitalization, 24h Trading Volume, Price Change (24h %)
# $null means that particular column's value did not change for that row
$data = @(
    @(2, 99002, 1957894470032, 118089762693, 2.41144),
    @(3, 3381.18, 406862764584, 57938597348, 9.056800000000001),
    @(4, 1.002, 130949157622, 173551742236, 0.11938),
    @(5, 262.44, 124220173708, 15298491107, 10.80211),
    @(6, 635.4, 92659978087, 2511013877, 5.36731),
    @(7, 1.42, 80591370672, 17643086232, 29.46279),
    @(8, 0.395552, 58058582415, 10368763281, 4.5807),
    @(9, 1.001, 38260938357, 10560509915, 0.04934),
    @(10, 3383.67, 33041881744, 148028738, 9.25051),
    @(11, 0.891699, 31872011634, 3682783069, 14.8535),
    @(12, 0.20108, 17356540142, 1103322051, 2.57305),
    @(13, 36.36, 14842131964, 1050149638, 9.352589999999999),
    @(14, 0.00002501, 14692282099, 1619515324, 5.795),
    @(15, 98749, 14410037964, 901639215, 2.73716),
    @(16, 4000.57, 14392682837, 167702447, 9.100989999999999),
    @(17, 5.56, 14148936409, 619899192, 4.88377),
    @(18, $null, 10308929159, 2225166820, 4.19795),
    @(19, 494.24, 9777605296, 2283872731, 6.40315),
    @(20, 3377.99, 9686429221, 2254650087, 9.13232),
    @(21, 15.15, 9490558072, 1228707366, 6.48594),
    @(22, 0.00002141, 8987945059, 7056904336, 15.03015),
    @(23, 6.23, 8921157612, 821158717, 11.8257),
    @(24, 0.287482, 8619859512, 2353187735, 22.54928),
    @(25, 8.73, 8070751002, 3455613, 2.38203),
    @(26, 5.8, 7040976947, 1021697502, 7.08168),
    @(27, 90.36, 6794019549, 1430474849, 6.60149),
    @(28, $null, 6439611393, 899442878, 6.00717),
    @(29, 3551.03, 6104576415, 99949949, 8.69844),
    @(30, 9.369999999999999, 5606275386, 844571896, 8.29304),
    @(31, 0.195206, 5297479520, 115772728, 11.89512),
    @(32, 1, 5243513951, 15727350, -0.10343),
    @(33, 0.132078, 5038922708, 855228186, 5.7654),
    @(34, 9.699999999999999, 4588051994, 272484306, 9.55508),
    @(35, 27.95, 4176985543, 904902335, 7.9547),
    @(36, 0.00005227, 3918665466, 1811077711, 8.74653),
    @(37, $null, 3820237069, 445811950, 2.8475),
    @(38, 0.150729, 3790237689, 155893008, 1.53851),
    @(39, 510.19, 3755619000, 288147466, 6.01131),
    @(40, 0.470041, 3732725346, 477379599, 9.45964),
    @(41, 1.004, 3687081210, 234094476, 0.20564),
    @(42, 24.82, 3568902481, 42268168, 3.171),
    @(43, 3.87, 3509830666, 307919878, 7.07413),
    @(44, 1.001, 3437170730, 162604667, 0.09223000000000001),
    @(45, 3.39, 3376536344, 1277498925, 9.75545),
    @(46, 1.29, 3350856732, 496633305, 6.08043),
    @(47, 0.77889, 3189387025, 1680163823, 15.18927),
    @(48, 160.55, 2960163016, 83785159, -0.64177),
    @(49, 1.95, 2930298939, 407497415, 4.39563),
    @(50, 4.7, 2810230521, 589963526, 10.18223),
    @(51, 46.64, 2798490620, 20268743, 6.34523)
)

foreach ($row in $data) {
    $r = $row[0]
    if ($null -ne $row[1]) {
        $ws1.Cells.Item($r, 3).Value = $row[1]
    }
    $ws1.Cells.Item($r, 4).Value = $row[2]
    $ws1.Cells.Item($r, 5).Value = $row[3]
    $ws1.Cells.Item($r, 6).Value = $row[4]
}

# Update "Top 5 by Market Cap" sheet (mirrors Market Capitalization for the top 5 coins)
$ws2.Cells.Item(2, 2).Value = 1957894470032
$ws2.Cells.Item(3, 2).Value = 406862764584
$ws2.Cells.Item(4, 2).Value = 130949157622
$ws2.Cells.Item(5, 2).Value = 124220173708
$ws2.Cells.Item(6, 2).Value = 92659978087

# Update "Summary" sheet
# The Average Price cell holds a currency-formatted string ("$4356.86"). Excel normally
# auto-converts a plain "$"-prefixed numeric string into a number when assigned directly,
# so the cell is temporarily switched to Text format to force the value to stay a string,
# then restored to the Normal style (format/value match the source data exactly).
$avgCell = $ws3.Cells.Item(2, 2)
$avgCell.Style = "Normal"
$avgCell.NumberFormat = "@"
$avgCell.Value = "$4356.86"
$avgCell.Style = "Normal"

$ws3.Cells.Item(3, 2).Value = "XRP (29.46%)"
$ws3.Cells.Item(4, 2).Value = "Monero (-0.64%)"
